$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct the revised B/D values for existing rows 356-358 ---
$ws.Range("B356").Value = 2078216858000
$ws.Range("D356").Value = 474869038022.1187

$ws.Range("B357").Value = 2082183969000
$ws.Range("D357").Value = 480662981370.7611

$ws.Range("B358").Value = 2118202312000
$ws.Range("D358").Value = 476707546473.4213

# --- Copy the date-column formatting down to the new rows ---
$ws.Range("A358").Copy()
$ws.Range("A359:A361").PasteSpecial(-4122)

# --- Append the three new rows of data ---
$ws.Range("A359").Value = 44986
$ws.Range("B359").Value = 2121975670000
$ws.Range("C359").Value = 0.2319647413593134
$ws.Range("D359").Value = 492223537462.3057

$ws.Range("A360").Value = 45017
$ws.Range("B360").Value = 2135028350000
$ws.Range("C360").Value = 0.240610187435336
$ws.Range("D360").Value = 513709571473.2562

$ws.Range("A361").Value = 45047
$ws.Range("B361").Value = 2140971740000
$ws.Range("C361").Value = 0.2359826316783085
$ws.Range("D361").Value = 505232145554.0873
